# Automatic update of files.
# The "Förändrad" (Changed) date column (C) for rows 2-33 moves forward
# one day, from serial 45617 (2024-11-21) to serial 45618 (2024-11-22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("C2:C33")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45617) {
        $cell.Value2 = 45618
    }
}
